$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.063.08'
$ws.Range('E2').Value = '  +2.32%  '
$ws.Range('D3').Value = '1.912.11'
$ws.Range('E3').Value = '  +2.62%  '
$ws.Range('E4').Value = '  -0.65%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.005'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4815'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.82%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3816'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07363'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9344'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.86'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07784'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').Value = '1.916.25'
$ws.Range('E13').Value = '  +2.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.500'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.645'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '92.19'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.006'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008868'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('D20').Value = '28.102.67'
$ws.Range('E20').Value = '  +2.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.79'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.56%  '
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('D23').Value = '2.137.29'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.91'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.919'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.49'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.134'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '116.92'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.969'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08960'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.308'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.268'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7806'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.680'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.619'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02057'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.113'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05318'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5506'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.44%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.001'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.037'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.495'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.64'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '108.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.37%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4831'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.005'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.650'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.06'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06086'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.16%  '
